$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 317 previously had no symbol (column B) value - fill it in now.
$ws.Cells.Item(317, 2).Value2 = "ECONOMICS:CNCBBS"

# New rows 318-326 reuse the same formatting as row 317 (A:G).
# Copy that row's formats once and paste (format-only) into each new
# row before setting values, so the existing style (date format on
# column A, etc.) is reused instead of a brand-new style being created.
$ws.Range("A317:G317").Copy()

$newRows = @(
    @{ Row = 318; Date = 45230; Amount = 43325980000000; HasSymbol = $true  },
    @{ Row = 319; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 320; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 321; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 322; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 323; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 324; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 325; Date = 45257; Amount = 44065463000000; HasSymbol = $true  },
    @{ Row = 326; Date = 45257; Amount = 44065463000000; HasSymbol = $false }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Range("A$($r):G$($r)").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value2 = $entry.Date
    if ($entry.HasSymbol) {
        $ws.Cells.Item($r, 2).Value2 = "ECONOMICS:CNCBBS"
    }
    else {
        $ws.Cells.Item($r, 2).Clear()
    }
    $ws.Cells.Item($r, 3).Value2 = $entry.Amount
    $ws.Cells.Item($r, 4).Value2 = $entry.Amount
    $ws.Cells.Item($r, 5).Value2 = $entry.Amount
    $ws.Cells.Item($r, 6).Value2 = $entry.Amount
    $ws.Cells.Item($r, 7).Value2 = 0
}
